$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Update the letter date: "September 19, 2025" -> "September 21, 2025"
# ---------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "September 19, 2025") {
        $p.Range.Text = "September 21, 2025"
        break
    }
}

# ---------------------------------------------------------------------
# 2) Split the standalone mailing-address paragraph into two paragraphs:
#      "2061 Holly Branch Court, Santa Clara CA 95050"
#    becomes
#      "2061 Holly Branch Court"
#      "Santa Clara, CA 95050"
#    (Only the free-standing address line above the table is touched;
#    the identical text inside the table stays untouched.)
# ---------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Information(12) -eq $false -and `
        $p.Range.Text.TrimEnd([char]13, [char]7) -eq "2061 Holly Branch Court, Santa Clara CA 95050") {

        $snippet = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
            '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
            '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
            '<w:p><w:pPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr>' + `
            '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">2061 Holly Branch Court</w:t></w:r></w:p>' + `
            '<w:p><w:pPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr>' + `
            '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">Santa Clara, CA 95050</w:t></w:r></w:p>' + `
            '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

        $p.Range.InsertXML($snippet)
        break
    }
}

# ---------------------------------------------------------------------
# 3) Remove the now-redundant empty "No Spacing" paragraph that trails
#    "... Board of Directors".
# ---------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Rancho Palma Grande Homeowners Association Board of Directors") {
        $next = $p.Next()
        if ($next -ne $null -and $next.Range.Text.TrimEnd([char]13, [char]7) -eq "" -and $next.Style.NameLocal -eq "No Spacing") {
            $next.Range.Delete()
        }
        break
    }
}
